# BIS-1002: remove "Internal Assignment" column (column O) from the
# sample-type export sheet. The header cell O4 ("Internal Assignment")
# and the per-property values O5:O9 ("FALSE") are cleared so the column
# no longer carries any content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O4:O9").ClearContents() | Out-Null

# Match the resulting selection left behind by the edit (the editor's
# cursor/selection ends up anchored on the now-empty header cell).
$ws.Range("O4:O9").Select() | Out-Null
